$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts everything right by one column)
$ws.Columns("A:A").Insert()

# New header for the inserted column
$ws.Range("A1").Value = "Match ID"

# Apply bold style (no border) to the new column for rows 1-17; row 18 (totals) stays default style
$ws.Range("A1:A17").Font.Bold = $true
$ws.Range("A18").Font.Bold = $false

# Fill in the Match ID values for the data rows
$ws.Range("A4:A17").Value = 16
$ws.Range("A18").Value = 16

# Writing into the hidden totals row can cause the engine to stamp an
# explicit row height; auto-fit it back so no custom height is persisted
$ws.Rows("18:18").AutoFit()

# Update the selection shown when the sheet is viewed
$ws.Range("A1:A17").Select()
